# Update cryptocurrency price (D) and 1h volume/change (E) columns
# to the latest scraped values, preserving exact text formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.025.28"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.359.57"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").Value = "'239.96"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("D7").Value = "'74.26"
$ws.Range("E7").Value = "  +2.30%  "
$ws.Range("D9").Value = "'0.596"
$ws.Range("E9").Value = "  +10.82%  "
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("D11").Value = "'57.20"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "'32.25"
$ws.Range("E12").Value = "  +11.55%  "
$ws.Range("D13").Value = "'7.31"
$ws.Range("E13").Value = "  +10.06%  "
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "2.711.58"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").Value = "'16.67"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").Value = "'0.901"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "2.364.60"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "43.912.69"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("D21").Value = "'6.76"
$ws.Range("E21").Value = "  +5.50%  "
$ws.Range("D22").Value = "'77.05"
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("D23").Value = "'257.59"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").Value = "'1.99"
$ws.Range("E24").Value = "  +26.45%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "'2.50"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -2.58%  "
$ws.Range("D28").Value = "'10.78"
$ws.Range("E28").Value = "  +2.67%  "
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").Value = "'22.77"
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("D31").Value = "'175.27"
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("D32").Value = "'0.129"
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("E33").Value = "  +4.00%  "
$ws.Range("E34").Value = "  +7.51%  "
$ws.Range("E35").Value = "  +2.04%  "
$ws.Range("D36").Value = "'5.49"
$ws.Range("E36").Value = "  +4.84%  "
$ws.Range("E37").Value = "  -5.01%  "
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("D39").Value = "'6.34"
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("E40").Value = "  +4.61%  "
$ws.Range("E41").Value = "  +15.46%  "
$ws.Range("E42").Value = "  +14.70%  "
$ws.Range("D43").Value = "'9.12"
$ws.Range("E43").Value = "  +3.47%  "
$ws.Range("D44").Value = "'19.04"
$ws.Range("E44").Value = "  -2.28%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  +6.68%  "
$ws.Range("D47").Value = "'58.33"
$ws.Range("E47").Value = "  +11.37%  "
$ws.Range("D48").Value = "'2.52"
$ws.Range("E48").Value = "  +8.47%  "
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("D51").Value = "'100.32"
$ws.Range("E51").Value = "  +2.28%  "
